$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.039.52"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "1.916.50"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5037"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4006"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08316"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.25%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.912.63"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.412"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.274"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "

$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06492"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.951"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("D23").Value = "30.077.02"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.192"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("E26").Value = "  +3.04%  "

$ws.Range("D27").Value = "2.142.75"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  -4.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1037"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.988"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.817"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02442"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.354"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06363"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6667"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2156"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.197"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.681"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.31%  "

$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6096"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.197"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.633"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.129"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.45%  "
